# Apply the cryptos-list refresh described by the commit:
# "Updated cryptos list on Sun May 21 22:32:13 UTC 2023 with GitHub Actions"
#
# Column D ("Price") and column E ("Volume(1h)") values are refreshed for each
# coin row. Both columns hold plain TEXT in the workbook (prices sometimes use
# a '.' thousands separator e.g. 26.976.45, and volumes keep padding spaces +
# a trailing "%"), so a leading single-quote is used wherever the new value
# would otherwise be auto-recognised as a number, forcing Excel to store it as
# text exactly like the original cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.976.45'
$ws.Range("E2").Value = '  -1.53%  '
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("E4").Value = '  -0.50%  '
$ws.Range("D5").Value = '''309.62'
$ws.Range("E5").Value = '  -1.62%  '
$ws.Range("D6").Value = '''1.009'
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").Value = '''0.4629'
$ws.Range("E7").Value = '  -3.02%  '
$ws.Range("D8").Value = '''0.3642'
$ws.Range("E8").Value = '  -1.65%  '
$ws.Range("D9").Value = '''0.07287'
$ws.Range("E9").Value = '  -2.49%  '
$ws.Range("D10").Value = '''0.8655'
$ws.Range("E10").Value = '  -2.38%  '
$ws.Range("D11").Value = '''19.83'
$ws.Range("E11").Value = '  -3.22%  '
$ws.Range("D12").Value = '1.893.87'
$ws.Range("E12").Value = '  +2.62%  '
$ws.Range("D13").Value = '''0.07626'
$ws.Range("E13").Value = '  +3.46%  '
$ws.Range("D14").Value = '''93.27'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = '''5.331'
$ws.Range("E15").Value = '  -2.89%  '
$ws.Range("D16").Value = '''6.475'
$ws.Range("E16").Value = '  -1.80%  '
$ws.Range("E17").Value = '  -0.70%  '
$ws.Range("D18").Value = '''0.000008626'
$ws.Range("E18").Value = '  -2.66%  '
$ws.Range("D20").Value = '27.362.88'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '''14.48'
$ws.Range("E21").Value = '  -2.50%  '
$ws.Range("D22").Value = '''5.163'
$ws.Range("E22").Value = '  -3.54%  '
$ws.Range("E23").Value = '  -1.53%  '
$ws.Range("D24").Value = '2.102.64'
$ws.Range("E24").Value = '  +1.52%  '
$ws.Range("D25").Value = '''151.83'
$ws.Range("E25").Value = '  -0.49%  '
$ws.Range("D26").Value = '''1.861'
$ws.Range("E26").Value = '  -2.23%  '
$ws.Range("E27").Value = '  -2.30%  '
$ws.Range("D28").Value = '''2.103'
$ws.Range("E28").Value = '  -3.00%  '
$ws.Range("D29").Value = '''5.088'
$ws.Range("E29").Value = '  -3.54%  '
$ws.Range("D30").Value = '''115.89'
$ws.Range("E30").Value = '  -1.92%  '
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").Value = '''2.953'
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").Value = '''0.7318'
$ws.Range("E33").Value = '  -3.82%  '
$ws.Range("E34").Value = '  -3.30%  '
$ws.Range("D35").Value = '''4.432'
$ws.Range("E35").Value = '  -2.86%  '
$ws.Range("D37").Value = '''2.521'
$ws.Range("E37").Value = '  +5.73%  '
$ws.Range("D38").Value = '''1.075'
$ws.Range("E38").Value = '  -2.96%  '
$ws.Range("D39").Value = '''0.05271'
$ws.Range("E39").Value = '  -2.20%  '
$ws.Range("E40").Value = '  -2.58%  '
$ws.Range("D41").Value = '''2.928'
$ws.Range("E41").Value = '  -2.39%  '
$ws.Range("D42").Value = '''7.147'
$ws.Range("E42").Value = '  -2.31%  '
$ws.Range("D43").Value = '''0.5219'
$ws.Range("E43").Value = '  -2.73%  '
$ws.Range("E44").Value = '  -2.08%  '
$ws.Range("D45").Value = '''8.243'
$ws.Range("E45").Value = '  -3.79%  '
$ws.Range("D46").Value = '''0.4861'
$ws.Range("E46").Value = '  -2.69%  '
$ws.Range("D48").Value = '''10.15'
$ws.Range("E48").Value = '  -4.09%  '
$ws.Range("D49").Value = '''103.17'
$ws.Range("E49").Value = '  -1.96%  '
$ws.Range("D50").Value = '''1.635'
$ws.Range("E50").Value = '  -2.90%  '
$ws.Range("D51").Value = '''0.06222'
$ws.Range("E51").Value = '  -1.71%  '
